$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 60.9
$ws.Range("N2").Value = 54.83846622768671

$ws.Range("K3").Value = 58.5
$ws.Range("N3").Value = 54.83846622768671

$ws.Range("K4").Value = 51.7
$ws.Range("N4").Value = 54.83846622768671

$ws.Range("K5").Value = 48.9
$ws.Range("N5").Value = 54.83846622768671
